$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows: shift the "twelve months ended" labels forward by one year ---
# Row 8 and Row 24 have headers in E:I reading year-ending labels.
# Old: 1396/12, 1397/12, 1398/12, 1399/12, 1400/12 (E..I)
# New: 1397/12, 1398/12, 1399/12, 1400/12, 1401/12 (E..I)
$ws.Range("E8").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E24").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1401/12"

# --- Data rows: shift values one column to the left (drop oldest year, add new year) ---
$dataRows = @(10, 11, 12, 16, 17, 19, 20, 26, 27)
$newLastCol = @{
    10 = 734608
    11 = 1272203
    12 = 251440
    16 = 45810
    17 = 1053220
    19 = 407890
    20 = 3765171
    26 = 186
    27 = 297
}

foreach ($r in $dataRows) {
    $fVal = $ws.Cells.Item($r, 6).Value2
    $gVal = $ws.Cells.Item($r, 7).Value2
    $hVal = $ws.Cells.Item($r, 8).Value2
    $iVal = $ws.Cells.Item($r, 9).Value2

    $ws.Cells.Item($r, 5).Value = $fVal
    $ws.Cells.Item($r, 6).Value = $gVal
    $ws.Cells.Item($r, 7).Value = $hVal
    $ws.Cells.Item($r, 8).Value = $iVal
    $ws.Cells.Item($r, 9).Value = $newLastCol[$r]
}

Write-Output "done"
